$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.371.12"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.238.78"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.77"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.44"
$ws.Range("E7").Value = "  -3.69%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.42"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0958"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.51"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.252.28"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "42.272.90"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000110"
$ws.Range("E18").Value = "  +11.95%  "
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.09"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.32"
$ws.Range("E21").Value = "  +40.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.73"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("E23").Value = "  -5.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.79"
$ws.Range("E24").Value = "  +4.29%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.64"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.30"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("E28").Value = "  +6.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.95"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.90"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.90"
$ws.Range("E31").Value = "  +19.99%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.119"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.05"
$ws.Range("E34").Value = "  -9.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.57"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("E37").Value = "  +2.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.38"
$ws.Range("E38").Value = "  -6.41%  "
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  -4.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.66"
$ws.Range("E41").Value = "  +3.99%  "
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.85"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.86"
$ws.Range("E44").Value = "  -6.97%  "
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.14"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.13"
$ws.Range("E51").Value = "  -1.71%  "
